$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shared string value for A4: "191448213678" -> "0191448213678"
$ws.Range("A4").Value = "0191448213678"

# Update B4 value: 1 -> 2
$ws.Range("B4").Value = 2

# Update the selection to D8
$ws.Range("D8").Select()
